# Apply the upload changes: fill in the missing Hora inicio/fin + Multimedia
# entries for the last three rows (29-31) of Hoja1 and refresh the active
# selection/scroll position to match the author's final view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Row 29 (2020-12-01 / Jueves) -----------------------------------------
$ws.Range("B29").Value = 0.92361111111111116
$ws.Range("B29").NumberFormat = "h:mm"
$ws.Range("C29").Value = 0.93472222222222223
$ws.Range("C29").NumberFormat = "h:mm"
$ws.Range("D29").Value = 3
$ws.Range("F29").Value = "YOUTUBE"
$ws.Range("H29").Value = 16

# --- Row 30 (2020-12-02 / Viernes) -----------------------------------------
$ws.Range("B30").Value = 0.54305555555555551
$ws.Range("B30").NumberFormat = "h:mm"
$ws.Range("C30").Value = 0.55347222222222225
$ws.Range("C30").NumberFormat = "h:mm"
$ws.Range("D30").Value = 2
$ws.Range("F30").Value = "YOUTUBE"
$ws.Range("G30").Value = 15

# --- Row 31 (2020-12-02 / Viernes) -----------------------------------------
$ws.Range("B31").Value = 0.9159722222222223
$ws.Range("B31").NumberFormat = "h:mm"
$ws.Range("C31").Value = 0.92638888888888893
$ws.Range("C31").NumberFormat = "h:mm"
$ws.Range("D31").Value = 2
$ws.Range("F31").Value = "TV"
$ws.Range("H31").Value = 15

# --- Refresh the view: scroll position + active selection ------------------
$ws.Range("G32").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
